# Auto-generated edit script: update crypto price/volume table
# (values refreshed by the "Updated cryptos list" GitHub Action run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.904.20"
$ws.Range("E2").Value = "'  -0.55%  "

$ws.Range("D3").Value = "'2.576.47"
$ws.Range("E3").Value = "'  +0.32%  "

$ws.Range("E4").Value = "'  +0.03%  "

$ws.Range("D5").Value = "'581.66"
$ws.Range("E5").Value = "'  -0.50%  "

$ws.Range("D6").Value = "'144.03"
$ws.Range("E6").Value = "'  -2.88%  "

$ws.Range("E7").Value = "'  +0.07%  "

$ws.Range("D8").Value = "'0.590"
$ws.Range("E8").Value = "'  -2.17%  "

$ws.Range("D9").Value = "'0.106"
$ws.Range("E9").Value = "'  -2.17%  "

$ws.Range("D10").Value = "'5.59"
$ws.Range("E10").Value = "'  -1.32%  "

$ws.Range("E11").Value = "'  -0.64%  "

$ws.Range("D12").Value = "'0.349"
$ws.Range("E12").Value = "'  -2.17%  "

$ws.Range("D13").Value = "'27.00"
$ws.Range("E13").Value = "'  -1.93%  "

$ws.Range("D14").Value = "'3.042.62"
$ws.Range("E14").Value = "'  +0.48%  "

$ws.Range("D15").Value = "'62.905.00"
$ws.Range("E15").Value = "'  -0.44%  "

$ws.Range("D16").Value = "'0.0000145"
$ws.Range("E16").Value = "'  -2.30%  "

$ws.Range("D17").Value = "'2.569.46"
$ws.Range("E17").Value = "'  -0.53%  "

$ws.Range("D18").Value = "'11.07"
$ws.Range("E18").Value = "'  -2.45%  "

$ws.Range("D19").Value = "'341.33"
$ws.Range("E19").Value = "'  -0.59%  "

$ws.Range("D20").Value = "'4.34"
$ws.Range("E20").Value = "'  -1.89%  "

$ws.Range("D21").Value = "'6.62"
$ws.Range("E21").Value = "'  -3.73%  "

$ws.Range("E22").Value = "'  +0.02%  "

$ws.Range("D23").Value = "'67.46"
$ws.Range("E23").Value = "'  +1.24%  "

$ws.Range("D24").Value = "'1.58"
$ws.Range("E24").Value = "'  +6.40%  "

$ws.Range("D25").Value = "'1.60"
$ws.Range("E25").Value = "'  -2.43%  "

$ws.Range("D26").Value = "'0.165"
$ws.Range("E26").Value = "'  -3.52%  "

$ws.Range("D27").Value = "'7.96"
$ws.Range("E27").Value = "'  -2.67%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "'  -0.08%  "

$ws.Range("D29").Value = "'8.23"
$ws.Range("E29").Value = "'  -3.92%  "

$ws.Range("D30").Value = "'1.92"
$ws.Range("E30").Value = "'  -3.28%  "

$ws.Range("D31").Value = "'458.57"
$ws.Range("E31").Value = "'  -0.56%  "

$ws.Range("D32").Value = "'0.0₃0798"
$ws.Range("E32").Value = "'  -3.20%  "

$ws.Range("D33").Value = "'1.66"
$ws.Range("E33").Value = "'  +1.62%  "

$ws.Range("D34").Value = "'176.89"
$ws.Range("E34").Value = "'  +0.01%  "

$ws.Range("E35").Value = "'  +0.12%  "

$ws.Range("D36").Value = "'0.399"
$ws.Range("E36").Value = "'  -2.17%  "

$ws.Range("D37").Value = "'18.86"
$ws.Range("E37").Value = "'  -2.07%  "

$ws.Range("D38").Value = "'4.50"
$ws.Range("E38").Value = "'  -0.14%  "

$ws.Range("E39").Value = "'  +0.01%  "

$ws.Range("D40").Value = "'1.69"
$ws.Range("E40").Value = "'  -3.49%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'40.07"
$ws.Range("E41").Value = "'  +0.91%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'158.27"
$ws.Range("E42").Value = "'  +4.46%  "

$ws.Range("D43").Value = "'3.68"
$ws.Range("E43").Value = "'  -3.70%  "

$ws.Range("D44").Value = "'21.23"
$ws.Range("E44").Value = "'  +0.44%  "

$ws.Range("D45").Value = "'0.634"
$ws.Range("E45").Value = "'  +2.87%  "

$ws.Range("D46").Value = "'0.0537"
$ws.Range("E46").Value = "'  -2.82%  "

$ws.Range("D47").Value = "'0.0961"
$ws.Range("E47").Value = "'  -2.06%  "

$ws.Range("D48").Value = "'0.0236"
$ws.Range("E48").Value = "'  -1.68%  "

$ws.Range("D49").Value = "'17.97"
$ws.Range("E49").Value = "'  -2.39%  "

$ws.Range("E50").Value = "'  +0.11%  "

$ws.Range("D51").Value = "'1.67"
$ws.Range("E51").Value = "'  -4.57%  "
